$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range('D2').Value = "'63.803.35"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.25%  "
$ws.Range('E2').Style = 'Normal'

# Row 3: Ethereum -> Ethereum
$ws.Range('D3').Value = "'2.623.16"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.63%  "
$ws.Range('E3').Style = 'Normal'

# Row 4: TetherUSD -> TetherUSD
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.15%  "
$ws.Range('E4').Style = 'Normal'

# Row 5: BNB -> BNB
$ws.Range('D5').Value = "'598.50"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.74%  "
$ws.Range('E5').Style = 'Normal'

# Row 6: Solana -> Solana
$ws.Range('D6').Value = "'151.23"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +4.14%  "
$ws.Range('E6').Style = 'Normal'

# Row 7: USDC -> USDC
$ws.Range('D7').Value = "'1.00"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.14%  "
$ws.Range('E7').Style = 'Normal'

# Row 8: XRP -> XRP
$ws.Range('D8').Value = "'0.590"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.96%  "
$ws.Range('E8').Style = 'Normal'

# Row 9: Dogecoin -> Dogecoin
$ws.Range('E9').Value = "'  +1.80%  "
$ws.Range('E9').Style = 'Normal'

# Row 10: Toncoin -> Toncoin
$ws.Range('E10').Value = "'  +3.71%  "
$ws.Range('E10').Style = 'Normal'

# Row 11: Cardano -> Cardano
$ws.Range('D11').Value = "'0.385"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +4.44%  "
$ws.Range('E11').Style = 'Normal'

# Row 12: TRON -> TRON
$ws.Range('E12').Value = "'  -0.85%  "
$ws.Range('E12').Style = 'Normal'

# Row 13: Avalanche -> Avalanche
$ws.Range('D13').Value = "'27.92"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.79%  "
$ws.Range('E13').Style = 'Normal'

# Row 14: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range('D14').Value = "'3.097.77"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.87%  "
$ws.Range('E14').Style = 'Normal'

# Row 15: WrappedBTC -> WrappedBTC
$ws.Range('D15').Value = "'63.647.59"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.41%  "
$ws.Range('E15').Style = 'Normal'

# Row 16: ShibaInu -> ShibaInu
$ws.Range('E16').Value = "'  +5.17%  "
$ws.Range('E16').Style = 'Normal'

# Row 17: WrappedEther -> WrappedEther
$ws.Range('D17').Value = "'2.626.64"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.92%  "
$ws.Range('E17').Style = 'Normal'

# Row 18: Chainlink -> Chainlink
$ws.Range('D18').Value = "'12.41"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +8.62%  "
$ws.Range('E18').Style = 'Normal'

# Row 19: Polkadot -> Polkadot
$ws.Range('D19').Value = "'4.72"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +4.54%  "
$ws.Range('E19').Style = 'Normal'

# Row 20: BitcoinCash -> BitcoinCash
$ws.Range('D20').Value = "'348.83"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +2.12%  "
$ws.Range('E20').Style = 'Normal'

# Row 21: Uniswap -> Uniswap
$ws.Range('E21').Value = "'  +0.86%  "
$ws.Range('E21').Style = 'Normal'

# Row 22: Dai -> Dai
$ws.Range('D22').Value = "'0.998"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -0.15%  "
$ws.Range('E22').Style = 'Normal'

# Row 23: Litecoin -> LEO
$ws.Range('B23').Value = "'LEO"
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range('C23').Style = 'Normal'
$ws.Range('D23').Value = "'5.76"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.08%  "
$ws.Range('E23').Style = 'Normal'

# Row 24: SuiNetwork -> Litecoin
$ws.Range('B24').Value = "'Litecoin"
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').Value = "'67.12"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.67%  "
$ws.Range('E24').Style = 'Normal'

# Row 25: InternetComputer(DFINITY) -> SuiNetwork
$ws.Range('B25').Value = "'SuiNetwork"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'1.73"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +9.27%  "
$ws.Range('E25').Style = 'Normal'

# Row 26: Fetch.AI -> InternetComputer(DFINITY)
$ws.Range('B26').Value = "'InternetComputer(DFINITY)"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'9.36"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +4.32%  "
$ws.Range('E26').Style = 'Normal'

# Row 27: Bittensor -> Fetch.AI
$ws.Range('B27').Value = "'Fetch.AI"
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = "'1.68"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.56%  "
$ws.Range('E27').Style = 'Normal'

# Row 28: Aptos -> Bittensor
$ws.Range('B28').Value = "'Bittensor"
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = "'555.68"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.10%  "
$ws.Range('E28').Style = 'Normal'

# Row 29: Kaspa -> Aptos
$ws.Range('B29').Value = "'Aptos"
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = "'8.12"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +5.06%  "
$ws.Range('E29').Style = 'Normal'

# Row 30: Binance-PegBSC-USD -> Kaspa
$ws.Range('B30').Value = "'Kaspa"
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = "'0.162"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.25%  "
$ws.Range('E30').Style = 'Normal'

# Row 31: PancakeSwap -> Binance-PegBSC-USD
$ws.Range('B31').Value = "'Binance-PegBSC-USD"
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = "'1.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  -0.16%  "
$ws.Range('E31').Style = 'Normal'

# Row 32: PEPE -> PancakeSwap
$ws.Range('B32').Value = "'PancakeSwap"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'2.06"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.06%  "
$ws.Range('E32').Style = 'Normal'

# Row 33: ImmutableX -> PEPE
$ws.Range('B33').Value = "'PEPE"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'0.0₃0853"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.31%  "
$ws.Range('E33').Style = 'Normal'

# Row 34: NEARProtocol -> ImmutableX
$ws.Range('B34').Value = "'ImmutableX"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'1.75"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.29%  "
$ws.Range('E34').Style = 'Normal'

# Row 35: Monero -> NEARProtocol
$ws.Range('B35').Value = "'NEARProtocol"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'5.30"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.81%  "
$ws.Range('E35').Style = 'Normal'

# Row 36: PolygonEcosystemToken -> Monero
$ws.Range('B36').Value = "'Monero"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'167.62"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.11%  "
$ws.Range('E36').Style = 'Normal'

# Row 37: FirstDigitalUSD -> PolygonEcosystemToken
$ws.Range('B37').Value = "'PolygonEcosystemToken"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'0.418"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +4.22%  "
$ws.Range('E37').Style = 'Normal'

# Row 38: EthereumClassic -> FirstDigitalUSD
$ws.Range('B38').Value = "'FirstDigitalUSD"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'1.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.04%  "
$ws.Range('E38').Style = 'Normal'

# Row 39: Stacks -> EthereumClassic
$ws.Range('B39').Value = "'EthereumClassic"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'19.61"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +3.45%  "
$ws.Range('E39').Style = 'Normal'

# Row 40: USDe -> Stacks
$ws.Range('B40').Value = "'Stacks"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'1.94"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.60%  "
$ws.Range('E40').Style = 'Normal'

# Row 41: Aave -> USDe
$ws.Range('B41').Value = "'USDe"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'0.999"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.01%  "
$ws.Range('E41').Style = 'Normal'

# Row 42: OKB -> Aave
$ws.Range('B42').Value = "'Aave"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'167.97"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +1.96%  "
$ws.Range('E42').Style = 'Normal'

# Row 43: Filecoin -> OKB
$ws.Range('B43').Value = "'OKB"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'39.72"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.41%  "
$ws.Range('E43').Style = 'Normal'

# Row 44: Hedera -> Filecoin
$ws.Range('B44').Value = "'Filecoin"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'3.98"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +6.53%  "
$ws.Range('E44').Style = 'Normal'

# Row 45: InjectiveProtocol -> Hedera
$ws.Range('B45').Value = "'Hedera"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'0.0590"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +4.93%  "
$ws.Range('E45').Style = 'Normal'

# Row 46: Mantle -> InjectiveProtocol
$ws.Range('B46').Value = "'InjectiveProtocol"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'22.06"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +1.75%  "
$ws.Range('E46').Style = 'Normal'

# Row 47: VeChain -> Mantle
$ws.Range('B47').Value = "'Mantle"
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = "'0.634"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.88%  "
$ws.Range('E47').Style = 'Normal'

# Row 48: dogwifhat -> VeChain
$ws.Range('B48').Value = "'VeChain"
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = "'0.0253"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +4.05%  "
$ws.Range('E48').Style = 'Normal'

# Row 49: BabyDogeCoin -> dogwifhat
$ws.Range('B49').Value = "'dogwifhat"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'2.02"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +6.27%  "
$ws.Range('E49').Style = 'Normal'

# Row 50: Stellar -> Stellar
$ws.Range('D50').Value = "'0.0970"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.70%  "
$ws.Range('E50').Style = 'Normal'

# Row 51: EnergySwap -> EnergySwap
$ws.Range('D51').Value = "'19.47"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +4.76%  "
$ws.Range('E51').Style = 'Normal'

